# Config sheet gains an "ExcelTest" run flag, a CustomerNumber column, and
# several new TestCaseName rows (DashboardPageTest / UpdatePageTest /
# SubmitPageTest) used by the Jenkins run. TestData gains a duplicate
# LoginPageTest credentials row, and the active sheet/selection moves back
# to Config.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Config")
$ws2 = $wb.Worksheets.Item("TestData")

# --- Config: row 3 ("HomePageTest") execution flags flip from No -> Yes ---
$ws1.Range("A3").Value = "Yes"
$ws1.Range("C3").Value = "Yes"

# --- Config: new "CustomerNumber" header column (D), styled like the
#     other headers ---
$ws1.Range("D1").Value = "CustomerNumber"
$ws1.Range("C1").Copy() | Out-Null
$ws1.Range("D1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Config: new test-case rows ---
$ws1.Range("A4").Value = "Yes"
$ws1.Range("B4").Value = "ExcelTest"
$ws1.Range("C4").Value = "Yes"

$ws1.Range("A5").Value = "No"
$ws1.Range("B5").Value = "DashboardPageTest"
$ws1.Range("C5").Value = "Yes"

$ws1.Range("A6").Value = "No"
$ws1.Range("B6").Value = "UpdatePageTest"
$ws1.Range("C6").Value = "Yes"

$ws1.Range("A7").Value = "Yes"
$ws1.Range("B7").Value = "SubmitPageTest"
$ws1.Range("C7").Value = "Yes"

# --- Config: widen columns to fit the new / longer content ---
$ws1.Columns.Item(1).ColumnWidth = 9.33203125
$ws1.Columns.Item(2).ColumnWidth = 17.5
$ws1.Columns.Item(3).ColumnWidth = 11
$ws1.Columns.Item(4).ColumnWidth = 16

# --- TestData: duplicate the LoginPageTest credentials into row 3 ---
$ws2.Range("A3").Value = "LoginPageTest"
$ws2.Range("B3").Value = "madhur_b"
$ws2.Range("C3").Value = "madhurb"
$ws2.Range("B2").Copy() | Out-Null
$ws2.Range("B3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Selection: TestData highlights the new row, then focus returns to
#     Config (which becomes the active tab) at B3 ---
$ws2.Activate() | Out-Null
$ws2.Range("A3:C3").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("B3").Select() | Out-Null
